$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 86; this shifts the existing rows 86-102 down to 87-103,
# growing the sheet's used range from A1:T102 to A1:T103.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with a new weekly price record
# (same Mercado/Producto metadata as the surrounding rows, new Fecha/Calidad/prices).
$ws.Range("A86").Value = 1
$ws.Range("B86").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C86").Value = "Arica y Parinacota"
$ws.Range("D86").Value = 44641
$ws.Range("E86").Value = 15
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100108
$ws.Range("H86").Value = "Tropicales y subtropicales"
$ws.Range("I86").Value = 100108003
$ws.Range("J86").Value = "Maracuyá"
$ws.Range("K86").Value = "Sin especificar"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 130
$ws.Range("N86").Value = 15000
$ws.Range("O86").Value = 16000
$ws.Range("P86").Value = 15500
$ws.Range("Q86").Value = "$/caja 20 kilos"
$ws.Range("R86").Value = "Región de Arica y Parinacota"
$ws.Range("S86").Value = 775
$ws.Range("T86").Value = 20
